# Refresh the "cryptos" sheet (Price / Volume(1h) columns, plus a couple of
# row re-orderings) to match the latest scrape, mirroring the GitHub Actions
# "Updated cryptos list" commit.
#
# Note: column D ("Price") holds plain-text figures (e.g. "256.77", "1.00").
# When such a string is assigned straight to .Value, Excel's smart-typing
# will silently coerce it to a Number (dropping significant trailing zeros,
# e.g. "1.00" -> 1). Prefixing with a leading apostrophe forces Excel to
# keep it as Text, exactly like the original inline-string cells. Values
# that already contain multiple "." separators (e.g. "98.444.00") are never
# number-parseable, so no apostrophe is required for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.444.00"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.425.75"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'256.77"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'678.82"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  -4.97%  "
$ws.Range("D8").Value = "'0.438"
$ws.Range("E8").Value = "  -4.35%  "
$ws.Range("D9").Value = "'1.08"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "3.422.89"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").Value = "'0.217"
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("D13").Value = "'42.27"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "'6.43"
$ws.Range("E14").Value = "  +13.98%  "
$ws.Range("D15").Value = "98.114.12"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "'0.0000269"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "4.063.90"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "'9.07"
$ws.Range("E18").Value = "  +17.66%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.426.75"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("B20").Value = "Stellar"
$ws.Range("C20").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D20").Value = "'0.580"
$ws.Range("E20").Value = "  +34.10%  "
$ws.Range("D21").Value = "'17.83"
$ws.Range("D22").Value = "'11.16"
$ws.Range("E22").Value = "  +6.41%  "
$ws.Range("E23").Value = "  -4.21%  "
$ws.Range("D24").Value = "'512.70"
$ws.Range("E24").Value = "  -3.23%  "
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "'6.64"
$ws.Range("E26").Value = "  +3.47%  "
$ws.Range("D27").Value = "'101.63"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "'12.96"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "3.605.69"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").Value = "'0.151"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("E31").Value = "  +5.15%  "
$ws.Range("D32").Value = "'0.199"
$ws.Range("E32").Value = "  +4.01%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'2.58"
$ws.Range("E34").Value = "  +20.34%  "
$ws.Range("D35").Value = "'0.580"
$ws.Range("E35").Value = "  +6.89%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "'30.22"
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("D38").Value = "'8.10"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  +11.87%  "
$ws.Range("D40").Value = "'543.51"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'0.887"
$ws.Range("E43").Value = "  +7.20%  "
$ws.Range("D44").Value = "'24.73"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.84"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").Value = "'5.95"
$ws.Range("E46").Value = "  +15.55%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0438"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").Value = "'9.08"
$ws.Range("E48").Value = "  +13.23%  "
$ws.Range("E49").Value = "  +13.83%  "
$ws.Range("D50").Value = "'3.30"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").Value = "'54.21"
$ws.Range("E51").Value = "  +9.26%  "
